# Auto-generated edit script: updates column F ("想去人数") values across all 4 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 235
$ws.Cells.Item(3, 6).Value = 425
$ws.Cells.Item(4, 6).Value = 166
$ws.Cells.Item(5, 6).Value = 3905
$ws.Cells.Item(7, 6).Value = 2569
$ws.Cells.Item(8, 6).Value = 81
$ws.Cells.Item(9, 6).Value = 3179
$ws.Cells.Item(11, 6).Value = 2328
$ws.Cells.Item(14, 6).Value = 325
$ws.Cells.Item(15, 6).Value = 463
$ws.Cells.Item(16, 6).Value = 14
$ws.Cells.Item(17, 6).Value = 26
$ws.Cells.Item(18, 6).Value = 217
$ws.Cells.Item(19, 6).Value = 350
$ws.Cells.Item(20, 6).Value = 310
$ws.Cells.Item(21, 6).Value = 425
$ws.Cells.Item(23, 6).Value = 1417
$ws.Cells.Item(25, 6).Value = 13
$ws.Cells.Item(27, 6).Value = 137
$ws.Cells.Item(28, 6).Value = 156
$ws.Cells.Item(29, 6).Value = 34
$ws.Cells.Item(31, 6).Value = 64
$ws.Cells.Item(32, 6).Value = 4366
$ws.Cells.Item(33, 6).Value = 4170
$ws.Cells.Item(34, 6).Value = 82
$ws.Cells.Item(35, 6).Value = 131
$ws.Cells.Item(36, 6).Value = 62
$ws.Cells.Item(38, 6).Value = 1147
$ws.Cells.Item(40, 6).Value = 484
$ws.Cells.Item(42, 6).Value = 1322
$ws.Cells.Item(43, 6).Value = 179
$ws.Cells.Item(44, 6).Value = 133
$ws.Cells.Item(45, 6).Value = 110
$ws.Cells.Item(46, 6).Value = 42
$ws.Cells.Item(49, 6).Value = 2

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 12
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(8, 6).Value = 26
$ws.Cells.Item(15, 6).Value = 45
$ws.Cells.Item(16, 6).Value = 213

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 154
$ws.Cells.Item(4, 6).Value = 2323

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 235
$ws.Cells.Item(4, 6).Value = 154
$ws.Cells.Item(5, 6).Value = 425
$ws.Cells.Item(6, 6).Value = 12
$ws.Cells.Item(8, 6).Value = 166
$ws.Cells.Item(10, 6).Value = 3905
$ws.Cells.Item(12, 6).Value = 2569
$ws.Cells.Item(13, 6).Value = 81
$ws.Cells.Item(14, 6).Value = 3179
$ws.Cells.Item(17, 6).Value = 2328
$ws.Cells.Item(20, 6).Value = 326
$ws.Cells.Item(21, 6).Value = 14
$ws.Cells.Item(22, 6).Value = 26
$ws.Cells.Item(23, 6).Value = 350
$ws.Cells.Item(24, 6).Value = 310
$ws.Cells.Item(25, 6).Value = 425
$ws.Cells.Item(27, 6).Value = 1417
$ws.Cells.Item(30, 6).Value = 156
$ws.Cells.Item(32, 6).Value = 64
$ws.Cells.Item(34, 6).Value = 4366
$ws.Cells.Item(35, 6).Value = 4170
$ws.Cells.Item(36, 6).Value = 82
$ws.Cells.Item(38, 6).Value = 1147
$ws.Cells.Item(44, 6).Value = 45
$ws.Cells.Item(45, 6).Value = 1322
$ws.Cells.Item(46, 6).Value = 179
$ws.Cells.Item(47, 6).Value = 110
$ws.Cells.Item(48, 6).Value = 42
$ws.Cells.Item(49, 6).Value = 213
